$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32 (G=5484)
$ws.Range("H32").Value = 2002
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2002
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2002
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2654
# Row 41 (G=5478)
$ws.Range("H41").Value = 870.5238000000001
$ws.Range("I41").Value = 1398.25
$ws.Range("J41").Value = 545.7692
$ws.Range("K41").Value = 1398.25
$ws.Range("L41").Value = 545.7692
$ws.Range("M41").Value = -958.25
$ws.Range("N41").Value = -1425.7692
# Row 43 (G=5472)
$ws.Range("H43").Value = 1580.4
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1580.4
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1580.4
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1718.4
# Row 69 (G=12616)
$ws.Range("H69").Value = 9580
$ws.Range("J69").Value = 9580
$ws.Range("L69").Value = 28740
$ws.Range("N69").Value = -30488
# Row 72 (G=12616)
$ws.Range("H72").Value = 9580
$ws.Range("J72").Value = 9580
$ws.Range("L72").Value = 86220
$ws.Range("N72").Value = -94956
# Row 88 (G=12608)
$ws.Range("H88").Value = 4445.4546
$ws.Range("J88").Value = 4790
$ws.Range("L88").Value = 4790
$ws.Range("N88").Value = -5602
# Row 91 (G=12608)
$ws.Range("H91").Value = 4445.4546
$ws.Range("J91").Value = 4790
$ws.Range("L91").Value = 4790
$ws.Range("N91").Value = -7598
# Row 135 (G=44047)
$ws.Range("H135").Value = 837.5833
$ws.Range("I135").Value = 798.64703
$ws.Range("J135").Value = 932.1429000000001
$ws.Range("K135").Value = 7187.82327
$ws.Range("L135").Value = 8389.286100000001
$ws.Range("M135").Value = -4652.82327
$ws.Range("N135").Value = -13459.2861
# Row 137 (G=44013)
$ws.Range("H137").Value = 1018.95
$ws.Range("I137").Value = 1045.625
$ws.Range("K137").Value = 3136.875
$ws.Range("M137").Value = -586.875
# Row 138 (G=44169)
$ws.Range("H138").Value = 4268.7144
$ws.Range("I138").Value = 2263
$ws.Range("J138").Value = 4937.2856
$ws.Range("K138").Value = 6789
$ws.Range("L138").Value = 14811.8568
$ws.Range("M138").Value = -1649
$ws.Range("N138").Value = -25091.8568

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61 (G=43999)
$ws.Range("H61").Value = 2520
$ws.Range("I61").Value = 1613.3334
$ws.Range("K61").Value = 1613.3334
$ws.Range("M61").Value = -1401.3334
# Row 74 (G=44000)
$ws.Range("H74").Value = 2045.6875
$ws.Range("I74").Value = 1262.2273
$ws.Range("J74").Value = 3769.3
$ws.Range("K74").Value = 1262.2273
$ws.Range("L74").Value = 3769.3
$ws.Range("M74").Value = -388.2273
$ws.Range("N74").Value = -5517.3
# Row 77 (G=44000)
$ws.Range("H77").Value = 2045.6875
$ws.Range("I77").Value = 1262.2273
$ws.Range("J77").Value = 3769.3
$ws.Range("K77").Value = 6311.136500000001
$ws.Range("L77").Value = 18846.5
$ws.Range("M77").Value = -1943.136500000001
$ws.Range("N77").Value = -27582.5
# Row 122 (G=36168)
$ws.Range("H122").Value = 2770.3572
$ws.Range("I122").Value = 3101.4443
$ws.Range("J122").Value = 2174.4
$ws.Range("K122").Value = 9304.332900000001
$ws.Range("L122").Value = 6523.200000000001
$ws.Range("M122").Value = -6854.332900000001
$ws.Range("N122").Value = -11423.2
# Row 136 (G=43999)
$ws.Range("H136").Value = 2520
$ws.Range("I136").Value = 1613.3334
$ws.Range("K136").Value = 4840.0002
$ws.Range("M136").Value = -2290.0002

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134 (G=43998)
$ws.Range("H134").Value = 1597.45
$ws.Range("I134").Value = 1635.7059
$ws.Range("J134").Value = 1380.6666
$ws.Range("K134").Value = 4907.1177
$ws.Range("L134").Value = 4141.9998
$ws.Range("M134").Value = -2372.1177
$ws.Range("N134").Value = -9211.9998

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22 (G=5367)
$ws.Range("H22").Value = 739.6
$ws.Range("I22").Value = 310.8889
$ws.Range("J22").Value = 1382.6666
$ws.Range("K22").Value = 310.8889
$ws.Range("L22").Value = 1382.6666
$ws.Range("M22").Value = 39.11110000000002
$ws.Range("N22").Value = -2082.6666
# Row 56 (G=1867)
$ws.Range("H56").Value = 6546.5
$ws.Range("I56").Value = 3093
$ws.Range("J56").Value = 10000
$ws.Range("K56").Value = 3093
$ws.Range("L56").Value = 10000
$ws.Range("M56").Value = -2248
$ws.Range("N56").Value = -11690
# Row 86 (G=12584)
$ws.Range("H86").Value = 1941
$ws.Range("J86").Value = 2269.2
$ws.Range("L86").Value = 2269.2
$ws.Range("N86").Value = -4515.2
# Row 87 (G=11929)
$ws.Range("H87").Value = 39999
$ws.Range("J87").Value = 39999
$ws.Range("L87").Value = 39999
$ws.Range("N87").Value = -42371
# Row 89 (G=12584)
$ws.Range("H89").Value = 1941
$ws.Range("J89").Value = 2269.2
$ws.Range("L89").Value = 11346
$ws.Range("N89").Value = -22578
# Row 90 (G=11929)
$ws.Range("H90").Value = 39999
$ws.Range("J90").Value = 39999
$ws.Range("L90").Value = 119997
$ws.Range("N90").Value = -131853
# Row 107 (G=27689)
$ws.Range("H107").Value = 871.2222
$ws.Range("I107").Value = 948.1818
$ws.Range("J107").Value = 750.2857
$ws.Range("K107").Value = 948.1818
$ws.Range("L107").Value = 750.2857
$ws.Range("M107").Value = 971.8182
$ws.Range("N107").Value = -4590.2857
# Row 134 (G=44020)
$ws.Range("H134").Value = 1469.7646
$ws.Range("I134").Value = 1515.8462
$ws.Range("J134").Value = 1320
$ws.Range("K134").Value = 4547.5386
$ws.Range("L134").Value = 3960
$ws.Range("M134").Value = -2012.5386
$ws.Range("N134").Value = -9030

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 117 (G=27870)
$ws.Range("H117").Value = 7520.1875
$ws.Range("I117").Value = 345
$ws.Range("J117").Value = 9911.916999999999
$ws.Range("K117").Value = 1035
$ws.Range("L117").Value = 29735.751
$ws.Range("M117").Value = 2407
$ws.Range("N117").Value = -36619.751
# Row 131 (G=36060)
$ws.Range("H131").Value = 802.47
$ws.Range("J131").Value = 814.7311999999999
$ws.Range("L131").Value = 2444.1936
$ws.Range("N131").Value = -12524.1936

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80 (G=12521)
$ws.Range("H80").Value = 2520.9092
$ws.Range("I80").Value = 2335
$ws.Range("J80").Value = 3016.6667
$ws.Range("K80").Value = 2335
$ws.Range("L80").Value = 3016.6667
$ws.Range("M80").Value = -1337
$ws.Range("N80").Value = -5012.6667
# Row 83 (G=12521)
$ws.Range("H83").Value = 2520.9092
$ws.Range("I83").Value = 2335
$ws.Range("J83").Value = 3016.6667
$ws.Range("K83").Value = 11675
$ws.Range("L83").Value = 15083.3335
$ws.Range("M83").Value = -6683
$ws.Range("N83").Value = -25067.3335
# Row 98 (G=18359)
$ws.Range("H98").Value = 14214
$ws.Range("J98").Value = 14214
$ws.Range("L98").Value = 14214
$ws.Range("N98").Value = -20204
# Row 128 (G=34544)
$ws.Range("H128").Value = 51500
$ws.Range("J128").Value = 51500
$ws.Range("L128").Value = 51500
$ws.Range("N128").Value = -61460
# Row 132 (G=44008)
$ws.Range("H132").Value = 2340.24
$ws.Range("I132").Value = 1952.5294
$ws.Range("J132").Value = 3164.125
$ws.Range("K132").Value = 5857.5882
$ws.Range("L132").Value = 9492.375
$ws.Range("M132").Value = -3327.5882
$ws.Range("N132").Value = -14552.375

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 20 (G=4308)
$ws.Range("H20").Value = 70006
$ws.Range("J20").Value = 70006
$ws.Range("L20").Value = 70006
$ws.Range("N20").Value = -70458
# Row 128 (G=34582)
$ws.Range("H128").Value = 47895
$ws.Range("J128").Value = 47895
$ws.Range("L128").Value = 47895
$ws.Range("N128").Value = -57855
# Row 132 (G=44058)
$ws.Range("H132").Value = 6690.923
$ws.Range("I132").Value = 8428.429
$ws.Range("J132").Value = 4663.8335
$ws.Range("K132").Value = 25285.287
$ws.Range("L132").Value = 13991.5005
$ws.Range("M132").Value = -22755.287
$ws.Range("N132").Value = -19051.5005

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 56 (G=10912)
$ws.Range("H56").Value = 37186
$ws.Range("J56").Value = 37186
$ws.Range("L56").Value = 37186
$ws.Range("N56").Value = -38614
# Row 132 (G=44029)
$ws.Range("H132").Value = 3934.3157
$ws.Range("I132").Value = 4404.5
$ws.Range("J132").Value = 3128.2856
$ws.Range("K132").Value = 13213.5
$ws.Range("L132").Value = 9384.856800000001
$ws.Range("M132").Value = -10683.5
$ws.Range("N132").Value = -14444.8568
# Row 136 (G=44031)
$ws.Range("H136").Value = 1531.6666
$ws.Range("I136").Value = 619.8823
$ws.Range("K136").Value = 1859.6469
$ws.Range("M136").Value = 690.3531
